$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# "Sprint No." value cell: "1" -> "2"
$t.Cell(2, 4).Range.Text = "2"

# "Review Date" value cell: "02/09/18" -> "02/21/18"
$t.Cell(3, 2).Range.Text = "02/21/18"
